# Update 13C-MFA fit results (val/WRES/SSRES/enrich_mfa columns) for SC WT-batch
# Rows correspond to measured metabolite mass-isotopomer distributions;
# row 11 (glc__D_d1) and row 44 (glc__D_d2) are substrate refs and are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "0.075491 0.188433 0.251416 0.250974 0.180762 0.057029"
$ws.Range("F2").Value = "0.11596774193548345 -2.102096774193551 -0.8756451612903253 1.0954838709677364 1.568709677419357 0.859677419354839"
$ws.Range("G2").Value = 9.598994042663902
$ws.Range("I2").Value = 0.490476

$ws.Range("E3").Value = "0.14791 0.110689 0.259274 0.235517 0.121841 0.125786"
$ws.Range("F3").Value = "0.013548387096774313 -0.2440322580645162 -0.13758064516128915 -0.04258064516128814 0.40612903225806557 0.1685483870967749"
$ws.Range("G3").Value = 0.2738261966701359
$ws.Range("I3").Value = 0.4904164

$ws.Range("E4").Value = "0.080694 0.124176 0.189436 0.247618 0.176211 0.12436 0.059684"
$ws.Range("F4").Value = "0.07016129032258142 -1.052741935483874 -0.5951612903225805 -0.16016129032257967 0.7461290322580676 1.0174193548387087 0.3258064516129034"
$ws.Range("G4").Value = 3.191057310093662
$ws.Range("I4").Value = 0.4884416666666667

$ws.Range("E5").Value = "0.267446 0.23548 0.007838 0.261383 0.232174"
$ws.Range("F5").Value = "-0.021129032258063264 0.472096774193546 0.7019354838709677 -0.17580645161290895 -0.2801612903225833"
$ws.Range("G5").Value = 0.8254334807492207
$ws.Range("I5").Value = 0.49100025

$ws.Range("E6").Value = "0.076919 0.124459 0.190274 0.248723 0.176997 0.124915 0.05995"
$ws.Range("F6").Value = "0.05435483870967806 -0.8182258064516142 -0.4169354838709684 -0.2053225806451607 0.5620967741935465 0.8849999999999987 0.29983870967742043"
$ws.Range("G6").Value = 2.077521514047866
$ws.Range("I6").Value = 0.4905731666666667

$ws.Range("E7").Value = "0.352091 0.004158 0.002535 0.322009 0.002481 0.002158 0.321609"
$ws.Range("F7").Value = "-0.025483870967740564 0.37725806451612903 0.4088709677419355 -0.18419354838710192 0.009354838709677472 0.3480645161290323 0.20177419354838744"
$ws.Range("G7").Value = 0.5060250520291383
$ws.Range("I7").Value = 0.4876038333333333

$ws.Range("E8").Value = "0.253865 0.22207 0.030885 0.246625 0.228716"
$ws.Range("F8").Value = "5.974032258064519 -2.9614516129032284 -2.7412903225806455 -3.0262903225806457 -0.12225806451613043"
$ws.Range("G8").Value = 61.14730985952139
$ws.Range("I8").Value = 0.4846447500000001

$ws.Range("E9").Value = "0.511454 0.003073 0.002194 0.484348"
$ws.Range("F9").Value = "0.00274193548386028 -0.08774193548387094 0.2609677419354839 -0.0035483870967777426"
$ws.Range("G9").Value = 0.07582291883454731
$ws.Range("I9").Value = 0.486835

$ws.Range("E10").Value = "0.487687 0.001384 0.00211 0.013232 0.00211 0.001384 0.487687"
$ws.Range("F10").Value = "0.01709677419354758 -0.06629032258064514 0.3403225806451613 0.025161290322580628 -0.36629032258064526 -0.645 -0.015645161290324797"
$ws.Range("G10").Value = 0.6715776274713842
$ws.Range("I10").Value = 0.497797

$ws.Range("E12").Value = "0.145124 0.111103 0.260701 0.236817 0.122514 0.12648"
$ws.Range("F12").Value = "-0.0014516129032272583 0.030483870967741085 -0.12548387096774208 -0.32919354838709475 0.5617741935483871 0.3056451612903239"
$ws.Range("G12").Value = 0.5340551768990629
$ws.Range("I12").Value = 0.4930824

$ws.Range("E13").Value = "0.143021 0.110998 0.260724 0.23684 0.122526 0.126492"
$ws.Range("F13").Value = "0.023870967741936974 -0.44629032258064605 0.4903225806451592 -0.7461290322580632 0.7653225806451612 0.00983870967742024"
$ws.Range("G13").Value = 1.582685093652442
$ws.Range("I13").Value = 0.493106

$ws.Range("E14").Value = "0.486799 0.026596 0.486519"
$ws.Range("F14").Value = "1.3885483870967708 -0.01370967741935456 -1.388709677419356"
$ws.Range("G14").Value = 3.856769146722158
$ws.Range("I14").Value = 0.499817

$ws.Range("E15").Value = "0.517991 0.003887 0.00254 0.47649"
$ws.Range("F15").Value = "0.0001612903225852832 -0.005645161290322631 0.15290322580645163 -0.0009677419354758854"
$ws.Range("G15").Value = 0.02341222684703433
$ws.Range("I15").Value = 0.479479

$ws.Range("E16").Value = "0.135089 0.118237 0.132936 0.242791 0.120099 0.130239 0.113497"
$ws.Range("F16").Value = "-0.027258064516129433 0.45967741935483736 -1.5996774193548404 1.242580645161293 -1.8270967741935469 0.6656451612903213 -0.06096774193548441"
$ws.Range("G16").Value = 8.1001040062435
$ws.Range("I16").Value = 0.4875091666666667

$ws.Range("E17").Value = "0.139764 0.007156 0.360784 0.011748 0.360749 0.005873 0.120221"
$ws.Range("F17").Value = "-0.00903225806451621 0.17774193548387093 -1.4529032258064543 -0.018225806451612946 1.1001612903225753 0.16467741935483862 1.0529032258064503"
$ws.Range("G17").Value = 4.489012460978142
$ws.Range("I17").Value = 0.4929425

$ws.Range("E18").Value = "0.13263 0.008309 0.366247 0.015228 0.366228 0.007613 0.122023"
$ws.Range("F18").Value = "-0.03887096774193407 0.6293548387096775 -1.3390322580645213 0.886774193548387 0.828225806451613 0.5962903225806452 1.3853225806451634"
$ws.Range("G18").Value = 5.937613111342372
$ws.Range("I18").Value = 0.5036003333333333

$ws.Range("E19").Value = "0.279114 0.215842 0.042598 0.238326 0.227099"
$ws.Range("F19").Value = "-0.04129032258065064 0.9359677419354844 0.5532258064516122 0.5035483870967764 -1.4709677419354867"
$ws.Range("G19").Value = 3.60110637356921
$ws.Range("I19").Value = 0.481103

$ws.Range("E20").Value = "0.144607 0.246039 0.120564 0.133582 0.244898 0.115202"
$ws.Range("F20").Value = "0.024354838709679395 -0.47790322580645056 1.5153225806451613 -0.2053225806451607 -0.2575806451612883 0.19016129032258056"
$ws.Range("G20").Value = 2.66985364203954
$ws.Range("I20").Value = 0.488703

$ws.Range("E21").Value = "0.144865 0.111082 0.260683 0.236801 0.122506 0.126471"
$ws.Range("F21").Value = "0.027258064516129433 -0.50516129032258 0.30935483870967345 -0.8661290322580623 0.8546774193548384 0.5683870967741934"
$ws.Range("G21").Value = 2.155348231009357
$ws.Range("I21").Value = 0.493046

$ws.Range("E22").Value = "0.154428 0.244174 0.113773 0.139607 0.242693 0.108307"
$ws.Range("F22").Value = "0.01370967741935512 -0.268709677419355 0.8037096774193551 -1.703870967741934 0.24338709677419074 1.3927419354838717"
$ws.Range("G22").Value = 5.620485744016645
$ws.Range("I22").Value = 0.4805696

$ws.Range("E23").Value = "0.520705 0.003847 0.002128 0.467726"
$ws.Range("F23").Value = "-0.019838709677416813 0.6204838709677419 -1.5267741935483872 0.023870967741932498"
$ws.Range("G23").Value = 2.717003069719043
$ws.Range("I23").Value = 0.470427

$ws.Range("E24").Value = "0.107328 0.030953 0.094331 0.15226 0.123569 0.123569 0.152258 0.094271 0.029587 0.093446"
$ws.Range("F24").Value = "0.01838709677419403 -0.1611290322580645 -0.5690322580645167 -1.1962903225806463 -0.36564516129032343 0.4612903225806454 0.6861290322580658 0.10629032258064403 0.10145161290322542 1.1720967741935477"
$ws.Range("G24").Value = 3.99386800208117
$ws.Range("I24").Value = 0.4932967777777778

$ws.Range("E25").Value = "0.151966 0.111476 0.260705 0.236814 0.122512 0.126478"
$ws.Range("F25").Value = "-0.0466129032258083 0.8832258064516143 -1.443225806451606 0.08758064516128838 1.1348387096774195 0.9891935483870989"
$ws.Range("G25").Value = 5.139194458896968
$ws.Range("I25").Value = 0.4931532000000001

$ws.Range("E26").Value = "0.49323 0.006357 0.006681 0.4929"
$ws.Range("F26").Value = "4.186612903225808 0.01709677419354842 -0.15032258064516124 -4.187580645161284"
$ws.Range("G26").Value = 35.08644843912587
$ws.Range("I26").Value = 0.499473

$ws.Range("E27").Value = "0.354675 0.013574 0.148944 0.148915 0.012252 0.330501"
$ws.Range("F27").Value = "-0.07161290322580645 1.3046774193548385 0.23677419354838664 -1.3474193548387097 0.813225806451613 0.49354838709677534"
$ws.Range("G27").Value = 4.483838735691988
$ws.Range("I27").Value = 0.491944

$ws.Range("E28").Value = "0.18386 0.024504 0.18476 0.116877 0.116876 0.184717 0.023192 0.166723"
$ws.Range("F28").Value = "0.11354838709677136 -1.4719354838709675 -0.294677419354838 -2.122258064516132 0.3062903225806448 -0.6946774193548396 2.2612903225806456 2.145806451612905"
$ws.Range("G28").Value = 17.06461110822062
$ws.Range("I28").Value = 0.4917081428571429

$ws.Range("E29").Value = "0.440452 0.06156 0.061528 0.44042"
$ws.Range("F29").Value = "0.014193548387102016 1.5188709677419348 -0.7762903225806448 -0.11806451612903394"
$ws.Range("G29").Value = 2.923736368366283
$ws.Range("I29").Value = 0.5019586666666666

$ws.Range("E30").Value = "0.266881 0.136667 0.215117 0.150305 0.233122"
$ws.Range("F30").Value = "-0.008387096774192982 0.18403225806451665 0.39354838709677387 0.6464516129032232 -0.8782258064516139"
$ws.Range("G30").Value = 1.377998803329863
$ws.Range("I30").Value = 0.487576

$ws.Range("E31").Value = "0.261164 0.237099 0.007898 0.263596 0.234139"
$ws.Range("F31").Value = "-0.10032258064516314 0.9956451612903245 0.6040322580645163 -0.3625806451612948 -0.5083870967741917"
$ws.Range("G31").Value = 1.756151040582732
$ws.Range("I31").Value = 0.49505975

$ws.Range("E32").Value = "0.051918 0.021808 0.090501 0.094869 0.115324 0.140226 0.140226 0.115322 0.094862 0.090425 0.020485 0.04095"
$ws.Range("F32").Value = "-0.1219354838709677 1.0900000000000003 -1.443870967741936 1.8285483870967743 -0.7996774193548394 -0.7572580645161299 -0.340645161290323 0.358870967741934 0.5187096774193544 0.23338709677419417 0.6996774193548384 1.4625806451612904"
$ws.Range("G32").Value = 11.04129198751301
$ws.Range("I32").Value = 0.5024067272727273

$ws.Range("E33").Value = "0.113102 0.031453 0.094044 0.151754 0.123158 0.123158 0.151752 0.093958 0.029488 0.093135"
$ws.Range("F33").Value = "-0.0077419354838719956 0.13177419354838737 -1.1908064516129033 -1.2993548387096787 -0.6924193548387105 0.8506451612903227 0.48161290322580463 0.6046774193548391 0.3630645161290322 1.56532258064516"
$ws.Range("G33").Value = 7.506445525494277
$ws.Range("I33").Value = 0.4917291111111111

$ws.Range("E34").Value = "0.055735 0.085961 0.106768 0.145457 0.117749 0.126536 0.141292 0.102655 0.084944 0.038236"
$ws.Range("F34").Value = "0.07999999999999943 -0.7575806451612914 -1.2462903225806448 -0.556612903225808 0.41064516129032363 -0.27338709677419387 0.6864516129032274 0.38483870967741923 0.9193548387096792 1.2127419354838707"
$ws.Range("G34").Value = 5.622028850156095
$ws.Range("I34").Value = 0.492173

$ws.Range("E35").Value = "0.234572 0.020744 0.24692 0.244227 0.021799 0.228368"
$ws.Range("F35").Value = "0.020967741935482457 -0.3864516129032262 0.1703225806451615 -0.41548387096774114 -0.22854838709677447 0.2956451612903228"
$ws.Range("G35").Value = 0.4910615504682622
$ws.Range("I35").Value = 0.4952602

$ws.Range("E36").Value = "0.431344 0.410762 0.140053 0.017396 0.000362 0.0"
$ws.Range("F36").Value = "0.0016129032258080647 -0.03161290322580226 0.09645161290322826 -0.0695161290322582 -0.010322580645161289 0.0"
$ws.Range("G36").Value = 0.01524393860561938
$ws.Range("I36").Value = 0.1489008

$ws.Range("E37").Value = "0.551235 0.026242 0.351087 0.012686 0.056325 0.001392"
$ws.Range("F37").Value = "0.010322580645171614 -0.18999999999999975 -0.015000000000001571 -0.037258064516129084 0.08919354838709645 -0.02387096774193547"
$ws.Range("G37").Value = 0.0463450312174819
$ws.Range("I37").Value = 0.1997468

$ws.Range("E38").Value = "0.543638 0.038738 0.339812 0.020353 0.054709 0.002644 3.2e-05"
$ws.Range("F38").Value = "0.009838709677415765 -0.14983870967741925 -0.029193548387094634 -0.1495161290322582 0.24177419354838717 0.059838709677419395 0.005161290322580645"
$ws.Range("G38").Value = 0.1078178459937564
$ws.Range("I38").Value = 0.1686115

$ws.Range("E39").Value = "0.732455 0.02299 0.237788 0.006067 0.000138"
$ws.Range("F39").Value = "0.18112903225806212 0.11225806451612882 -0.5730645161290324 0.17709677419354838 0.011935483870967743"
$ws.Range("G39").Value = 0.4053182622268463
$ws.Range("I39").Value = 0.12932975

$ws.Range("E40").Value = "0.540922 0.038587 0.34284 0.020535 0.055199 0.002667 3.2e-05"
$ws.Range("F40").Value = "0.022419354838709715 -0.3399999999999998 -0.07709677419355386 -0.014838709677419086 0.500483870967741 0.030161290322580664 0.005"
$ws.Range("G40").Value = 0.3736855359001039
$ws.Range("I40").Value = 0.1700325

$ws.Range("E41").Value = "0.519039 0.006102 0.422715 0.00225 0.046613 0.000106 0.0"
$ws.Range("F41").Value = "0.5850000000000076 0.052419354838709645 -0.765161290322576 -0.8617741935483872 0.4603225806451617 0.017096774193548388 0.0"
$ws.Range("G41").Value = 1.885288527575445
$ws.Range("I41").Value = 0.1742106666666667

$ws.Range("E42").Value = "0.692168 0.028174 0.266462 0.00797 0.013653"
$ws.Range("F42").Value = "0.21838709677418813 -0.393064516129032 -0.6375806451612968 -0.030645161290322732 2.2020967741935484"
$ws.Range("G42").Value = 5.458871045785646
$ws.Range("I42").Value = 0.159905

$ws.Range("E43").Value = "0.500391 0.005232 0.488692 0.001377 0.002084 7.4e-05 0.0"
$ws.Range("F43").Value = "-0.011451612903214876 0.1746774193548387 0.010645161290324274 -0.6998387096774193 0.27 0.0020967741935483866 -0.09290322580645162"
$ws.Range("G43").Value = 0.6020662851196668
$ws.Range("I43").Value = 0.1659088333333333

$ws.Range("E45").Value = "0.553937 0.02636 0.35384 0.012786 0.056767 0.001403"
$ws.Range("F45").Value = "-0.019193548387093588 0.35064516129032297 0.04032258064515685 0.6550000000000001 -0.3780645161290321 0.17274193548387096"
$ws.Range("G45").Value = 0.7267438865764826
$ws.Range("I45").Value = 0.2012962

$ws.Range("E46").Value = "0.549143 0.026112 0.352535 0.012739 0.056558 0.001398"
$ws.Range("F46").Value = "0.31967741935484506 0.014838709677419086 -0.48258064516128946 -0.011774193548387121 -0.09967741935483879 0.01516129032258063"
$ws.Range("G46").Value = 0.3456020031217514
$ws.Range("I46").Value = 0.2005242

$ws.Range("E47").Value = "0.749772 0.244148 0.000358"
$ws.Range("F47").Value = "0.34322580645161593 -1.053548387096776 -0.21258064516129033"
$ws.Range("G47").Value = 1.27295868886577
$ws.Range("I47").Value = 0.122432

$ws.Range("E48").Value = "0.71764 0.005353 0.277596 0.000221"
$ws.Range("F48").Value = "-0.0037096774193719788 0.11274193548387097 0.007419354838708144 0.014193548387096775"
$ws.Range("G48").Value = 0.01298100936524464
$ws.Range("I48").Value = 0.1870693333333333

$ws.Range("E49").Value = "0.543801 0.025221 0.350467 0.012602 0.056551 0.00144 3.3e-05"
$ws.Range("F49").Value = "1.0588709677419246 -0.24306451612903193 -1.4812903225806493 -0.029193548387096594 -0.8879032258064521 0.01903225806451614 -0.030806451612903223"
$ws.Range("G49").Value = 4.165044771071789
$ws.Range("I49").Value = 0.1662605

$ws.Range("E50").Value = "0.430669 0.016395 0.399801 0.009888 0.128026 0.001584 0.01366"
$ws.Range("F50").Value = "0.023870967741941453 -0.3579032258064514 -0.20112903225806422 -0.14645161290322595 0.5919354838709688 -0.036290322580645185 0.12967741935483887"
$ws.Range("G50").Value = 0.5590863423517183
$ws.Range("I50").Value = 0.2412741666666667

$ws.Range("E51").Value = "0.432494 0.019169 0.402234 0.0117 0.12851 0.001869 0.013678"
$ws.Range("F51").Value = "-0.027741935483871855 0.42548387096774165 -0.1432258064516129 0.3738709677419356 0.4040322580645173 0.17306451612903226 0.3516129032258064"
$ws.Range("G51").Value = 0.6589242976066606
$ws.Range("I51").Value = 0.2440316666666667

$ws.Range("E52").Value = "0.728907 0.025986 0.236498 0.00667 0.003573"
$ws.Range("F52").Value = "0.001451612903213828 -0.033870967741935785 -0.0075806451612889505 0.16451612903225807 0.13903225806451613"
$ws.Range("G52").Value = 0.04760234131113421
$ws.Range("I52").Value = 0.133321

$ws.Range("E53").Value = "0.559941 0.195893 0.18283 0.062306 0.001577 3.4e-05"
$ws.Range("F53").Value = "0.05129032258065169 -0.9529032258064511 0.6962903225806432 0.48887096774193534 0.12725806451612906 0.005483870967741935"
$ws.Range("G53").Value = 1.650694979188343
$ws.Range("I53").Value = 0.1509898

$ws.Range("E54").Value = "0.558832 0.026714 0.346014 0.012503 0.055508 0.001372"
$ws.Range("F54").Value = "-0.009516129032263105 0.17258064516129057 0.009354838709677822 -0.0774193548387096 -0.03129032258064512 0.08838709677419355"
$ws.Range("G54").Value = 0.04474726847034357
$ws.Range("I54").Value = 0.1970286

$ws.Range("E55").Value = "0.713601 0.039431 0.232202 0.011441 0.000276 3e-06"
$ws.Range("F55").Value = "0.0867741935483933 -0.4317741935483869 -0.19354838709677527 0.0022580645161290524 0.044516129032258066 0.00048387096774193554"
$ws.Range("G55").Value = 0.2334067117585862
$ws.Range("I55").Value = 0.1078554

$ws.Range("E56").Value = "0.776047 0.006388 0.223096 0.000178"
$ws.Range("F56").Value = "-0.0277419354838629 0.8474193548387097 0.072419354838706 0.02870967741935484"
$ws.Range("G56").Value = 0.7249579864724235
$ws.Range("I56").Value = 0.151038

$ws.Range("E57").Value = "0.491803 0.024557 0.358649 0.018739 0.086201 0.005531 0.006084 0.000568 1e-06 0.0"
$ws.Range("F57").Value = "0.14854838709677504 -1.4969354838709679 -0.2745161290322551 -1.0845161290322582 0.9250000000000007 0.2275806451612903 0.23435483870967747 0.05161290322580646 0.00016129032258064516 -0.00016129032258064516"
$ws.Range("G57").Value = 4.479420889698232
$ws.Range("I57").Value = 0.1345576666666667

$ws.Range("E58").Value = "0.562455 0.026775 0.358348 0.012949 0.05749 0.001421"
$ws.Range("F58").Value = "-0.022419354838709715 0.4188709677419353 -0.4198387096774229 0.5816129032258065 2.5082258064516134 0.06870967741935484"
$ws.Range("G58").Value = 6.986411342351722
$ws.Range("I58").Value = 0.2038766

$ws.Range("E59").Value = "0.73796 0.005425 0.239043 0.000366"
$ws.Range("F59").Value = "1.0925806451612818 -0.537258064516129 -3.3608064516129024 0.030322580645161294"
$ws.Range("G59").Value = 12.77831815816855
$ws.Range("I59").Value = 0.1615363333333333

$ws.Range("E60").Value = "0.465391 0.368747 0.116814 0.059941 0.011639 9e-06"
$ws.Range("F60").Value = "0.047903225806450275 -0.903225806451612 1.2222580645161294 2.6667741935483873 0.7020967741935484 -0.10016129032258066"
$ws.Range("G60").Value = 9.926683116545266
$ws.Range("I60").Value = 0.1657598

$ws.Range("E61").Value = "0.292456 0.294241 0.287798 0.088498 0.041497 0.008367 0.000725 1e-06"
$ws.Range("F61").Value = "0.008709677419354595 -0.10870967741935612 -0.14306451612903656 0.4064516129032249 0.5609677419354842 1.349516129032258 0.11693548387096774 0.0"
$ws.Range("G61").Value = 2.347116519250781
$ws.Range("I61").Value = 0.1925015714285714

$ws.Range("E62").Value = "0.70349 0.08784 0.199464 0.000168"
$ws.Range("F62").Value = "0.3441935483870918 -0.9333870967741936 -0.8025806451612917 -0.06596774193548388"
$ws.Range("G62").Value = 1.638168106139437
$ws.Range("I62").Value = 0.162424

$ws.Range("E63").Value = "0.679992 0.031095 0.265926 0.008298 0.020934"
$ws.Range("F63").Value = "-0.027580645161277617 0.6308064516129035 -0.012741935483870282 0.49112903225806454 -0.07435483870967736"
$ws.Range("G63").Value = 0.645576196670135
$ws.Range("I63").Value = 0.16789425

$ws.Range("E64").Value = "0.733847 0.023034 0.23824 0.006078 0.000138"
$ws.Range("F64").Value = "0.009516129032263105 0.029516129032257928 -0.03677419354838806 0.19129032258064524 0.02209677419354839"
$ws.Range("G64").Value = 0.03939435483870986
$ws.Range("I64").Value = 0.129575

$ws.Range("E65").Value = "0.245375 0.236952 0.230651 0.170269 0.074899 0.039205 0.009657 0.002937 0.000471 2.4e-05 0.0 0.0"
$ws.Range("F65").Value = "-0.16677419354838377 1.4706451612903206 -1.5512903225806478 0.22580645161290522 -0.4790322580645162 1.0025806451612902 0.7572580645161291 0.3446774193548387 0.07596774193548388 0.0038709677419354843 0.0 0.0"
$ws.Range("G65").Value = 6.580769276795008
$ws.Range("I65").Value = 0.1624697272727272

$ws.Range("E66").Value = "0.497261 0.025175 0.35746 0.018676 0.085907 0.005512 0.006063 0.000566 1e-06 0.0"
$ws.Range("F66").Value = "0.11145161290322536 -1.1204838709677416 -0.2533870967741918 -0.9206451612903223 0.889354838709676 0.40274193548387105 0.3104838709677419 0.03516129032258065 0.00016129032258064516 0.0"
$ws.Range("G66").Value = 3.230487773152961
$ws.Range("I66").Value = 0.1341843333333333

$ws.Range("E67").Value = "0.32401 0.285741 0.202977 0.139739 0.04256 0.015943 0.003428 0.000141 2e-06 0.0"
$ws.Range("F67").Value = "0.08548387096774684 -0.8762903225806397 0.13064516129031895 1.040483870967744 0.8295161290322584 0.9306451612903225 0.1816129032258065 0.02274193548387097 0.0003225806451612903 0.0"
$ws.Range("G67").Value = 3.46256503642039
$ws.Range("I67").Value = 0.1536042222222222

$ws.Range("E68").Value = "0.554577 0.016047 0.356635 0.006422 0.057315 0.000403"
$ws.Range("F68").Value = "0.521129032258062 -0.5796774193548387 -0.7183870967741958 -0.24387096774193553 -0.38403225806451635 0.017580645161290324"
$ws.Range("G68").Value = 1.330944302809574
$ws.Range("I68").Value = 0.1959716

